$wb = $excel.ActiveWorkbook

$wsCompleteness = $wb.Worksheets.Item("Completeness")
$wsInstructions = $wb.Worksheets.Item("Instructions")

# Add the "Template updated" note to the Instructions sheet (new shared string +
# a new red-font character style, matching the author's update note).
$wsInstructions.Activate()
$cell = $wsInstructions.Range("C1")
$cell.Value = "Template updated 12/8/22."
$cell.Font.Color = 255

# Leave the same cell selected/active, as the author's saved view shows.
$wsInstructions.Range("C3").Select()

# Restore the Completeness tab as the active sheet (it was active before the edit).
$wsCompleteness.Activate()
